# Update the NATMI LR-pairs TPM values for Tnc-Cntn1 sheet.
# Only numeric data values change; no structural/formula changes are involved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> FAPs via Tnc-Cntn1)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1346003333333333
$ws.Range("H2").Value = 0.403801
$ws.Range("I2").Value = 0.009651054304565105
$ws.Range("J2").Value = 0.009651054304565105
$ws.Range("Q2").Value = 0.004617240101111112
$ws.Range("R2").Value = 0.04155516091
$ws.Range("S2").Value = 0.009651054304565105
$ws.Range("T2").Value = 0.009651054304565105

# Row 3 (FAPs -> FAPs via Tnc-Cntn1)
$ws.Range("I3").Value = 0.8124788779145131
$ws.Range("J3").Value = 0.8124788779145132
$ws.Range("S3").Value = 0.8124788779145131
$ws.Range("T3").Value = 0.8124788779145132

# Row 4 (MuSCs -> FAPs via Tnc-Cntn1)
$ws.Range("G4").Value = 2.4807
$ws.Range("H4").Value = 7.4421
$ws.Range("I4").Value = 0.1778700677809217
$ws.Range("J4").Value = 0.1778700677809217
$ws.Range("Q4").Value = 0.085096279
$ws.Range("R4").Value = 0.765866511
$ws.Range("S4").Value = 0.1778700677809217
$ws.Range("T4").Value = 0.1778700677809217
